# Applies the scheduled-runner profit/price updates to the Ultros_Profits workbook.
# Each sheet (job) has a leve-profitability table; columns H-N hold price/profit
# figures that get refreshed by the scraper. We only touch the cells that changed.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1734.6666
$ws.Range("J17").Value = 1857.75
$ws.Range("L17").Value = 5573.25
$ws.Range("N17").Value = -5909.25
$ws.Range("H33").Value = 975
$ws.Range("I33").Value = 936
$ws.Range("J33").Value = 1170
$ws.Range("K33").Value = 936
$ws.Range("L33").Value = 1170
$ws.Range("M33").Value = -707
$ws.Range("N33").Value = -1628
$ws.Range("H96").Value = 191.38461
$ws.Range("I96").Value = 219.77777
$ws.Range("J96").Value = 127.5
$ws.Range("K96").Value = 659.33331
$ws.Range("L96").Value = 382.5
$ws.Range("M96").Value = 713.66669
$ws.Range("N96").Value = -3128.5
$ws.Range("H98").Value = 1390003.5
$ws.Range("I98").Value = 1587861.1
$ws.Range("K98").Value = 1587861.1
$ws.Range("M98").Value = -1586363.1
$ws.Range("H100").Value = 2962.0667
$ws.Range("I100").Value = 1720.2727
$ws.Range("J100").Value = 6377
$ws.Range("K100").Value = 1720.2727
$ws.Range("L100").Value = 6377
$ws.Range("M100").Value = -1179.2727
$ws.Range("N100").Value = -7459
$ws.Range("H101").Value = 1025.5333
$ws.Range("I101").Value = 406.8
$ws.Range("K101").Value = 1220.4
$ws.Range("M101").Value = 401.5999999999999
$ws.Range("H111").Value = 4074.25
$ws.Range("J111").Value = 5149.75
$ws.Range("L111").Value = 15449.25
$ws.Range("N111").Value = -21583.25
$ws.Range("H122").Value = 1390003.5
$ws.Range("I122").Value = 1587861.1
$ws.Range("K122").Value = 4763583.300000001
$ws.Range("M122").Value = -4761133.300000001
$ws.Range("H138").Value = 3261.3333
$ws.Range("I138").Value = 1210.2222
$ws.Range("K138").Value = 3630.6666
$ws.Range("M138").Value = 1509.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4774.3125
$ws.Range("I61").Value = 3217.2727
$ws.Range("J61").Value = 8199.799999999999
$ws.Range("K61").Value = 3217.2727
$ws.Range("L61").Value = 8199.799999999999
$ws.Range("M61").Value = -3005.2727
$ws.Range("N61").Value = -8623.799999999999
$ws.Range("H63").Value = 4509.5386
$ws.Range("I63").Value = 5493.125
$ws.Range("J63").Value = 2935.8
$ws.Range("K63").Value = 5493.125
$ws.Range("L63").Value = 2935.8
$ws.Range("M63").Value = -4807.125
$ws.Range("N63").Value = -4307.8
$ws.Range("H66").Value = 4509.5386
$ws.Range("I66").Value = 5493.125
$ws.Range("J66").Value = 2935.8
$ws.Range("K66").Value = 27465.625
$ws.Range("L66").Value = 14679
$ws.Range("M66").Value = -24033.625
$ws.Range("N66").Value = -21543
$ws.Range("H88").Value = 6412182.5
$ws.Range("I88").Value = 2861.2
$ws.Range("J88").Value = 10418008
$ws.Range("K88").Value = 2861.2
$ws.Range("L88").Value = 10418008
$ws.Range("M88").Value = -2455.2
$ws.Range("N88").Value = -10418820
$ws.Range("H91").Value = 6412182.5
$ws.Range("I91").Value = 2861.2
$ws.Range("J91").Value = 10418008
$ws.Range("K91").Value = 2861.2
$ws.Range("L91").Value = 10418008
$ws.Range("M91").Value = -1457.2
$ws.Range("N91").Value = -10420816
$ws.Range("H122").Value = 4626.6484
$ws.Range("I122").Value = 3689.5293
$ws.Range("J122").Value = 5423.2
$ws.Range("K122").Value = 11068.5879
$ws.Range("L122").Value = 16269.6
$ws.Range("M122").Value = -8618.5879
$ws.Range("N122").Value = -21169.6
$ws.Range("H132").Value = 4971.778
$ws.Range("I132").Value = 4451.273
$ws.Range("J132").Value = 5789.7144
$ws.Range("K132").Value = 13353.819
$ws.Range("L132").Value = 17369.1432
$ws.Range("M132").Value = -10823.819
$ws.Range("N132").Value = -22429.1432
$ws.Range("H136").Value = 4774.3125
$ws.Range("I136").Value = 3217.2727
$ws.Range("J136").Value = 8199.799999999999
$ws.Range("K136").Value = 9651.8181
$ws.Range("L136").Value = 24599.4
$ws.Range("M136").Value = -7101.8181
$ws.Range("N136").Value = -29699.4
$ws.Range("H138").Value = 82000
$ws.Range("J138").Value = 82000
$ws.Range("L138").Value = 82000
$ws.Range("N138").Value = -92280
$ws.Range("H139").Value = 79857.336
$ws.Range("J139").Value = 79857.336
$ws.Range("L139").Value = 79857.336
$ws.Range("N139").Value = -90137.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3511.4194
$ws.Range("I86").Value = 1444.421
$ws.Range("J86").Value = 6784.1665
$ws.Range("K86").Value = 1444.421
$ws.Range("L86").Value = 6784.1665
$ws.Range("M86").Value = -321.421
$ws.Range("N86").Value = -9030.166499999999
$ws.Range("H89").Value = 3511.4194
$ws.Range("I89").Value = 1444.421
$ws.Range("J89").Value = 6784.1665
$ws.Range("K89").Value = 7222.105
$ws.Range("L89").Value = 33920.8325
$ws.Range("M89").Value = -1606.105
$ws.Range("N89").Value = -45152.8325
$ws.Range("H134").Value = 5099.8423
$ws.Range("I134").Value = 3248.3
$ws.Range("J134").Value = 7157.1113
$ws.Range("K134").Value = 9744.900000000001
$ws.Range("L134").Value = 21471.3339
$ws.Range("M134").Value = -7209.900000000001
$ws.Range("N134").Value = -26541.3339

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H137").Value = 50000
$ws.Range("J137").Value = 50000
$ws.Range("L137").Value = 50000
$ws.Range("N137").Value = -60200

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 315.48
$ws.Range("I2").Value = 330.2857
$ws.Range("K2").Value = 330.2857
$ws.Range("M2").Value = -217.2857
$ws.Range("H80").Value = 18582678
$ws.Range("J80").Value = 27780270
$ws.Range("L80").Value = 27780270
$ws.Range("N80").Value = -27782266
$ws.Range("H83").Value = 18582678
$ws.Range("J83").Value = 27780270
$ws.Range("L83").Value = 138901350
$ws.Range("N83").Value = -138911334
$ws.Range("H102").Value = 1914.5405
$ws.Range("I102").Value = 1489.9062
$ws.Range("K102").Value = 1489.9062
$ws.Range("M102").Value = 132.0938000000001
$ws.Range("H122").Value = 5426.2383
$ws.Range("I122").Value = 4378.3
$ws.Range("J122").Value = 6378.909
$ws.Range("K122").Value = 13134.9
$ws.Range("L122").Value = 19136.727
$ws.Range("M122").Value = -10684.9
$ws.Range("N122").Value = -24036.727
$ws.Range("H126").Value = 4667
$ws.Range("I126").Value = 3399.4614
$ws.Range("K126").Value = 10198.3842
$ws.Range("M126").Value = -7728.3842

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7498.353
$ws.Range("I7").Value = 10588.143
$ws.Range("J7").Value = 5335.5
$ws.Range("K7").Value = 10588.143
$ws.Range("L7").Value = 5335.5
$ws.Range("M7").Value = -10476.143
$ws.Range("N7").Value = -5559.5
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").Value = $null
$ws.Range("H118").Value = 49636.363
$ws.Range("J118").Value = 49636.363
$ws.Range("L118").Value = 49636.363
$ws.Range("N118").Value = -52950.363
$ws.Range("H126").Value = 7498.353
$ws.Range("I126").Value = 10588.143
$ws.Range("J126").Value = 5335.5
$ws.Range("K126").Value = 31764.429
$ws.Range("L126").Value = 16006.5
$ws.Range("M126").Value = -29294.429
$ws.Range("N126").Value = -20946.5
$ws.Range("H134").Value = 63525
$ws.Range("J134").Value = 63525
$ws.Range("L134").Value = 63525
$ws.Range("N134").Value = -73665
$ws.Range("H135").Value = 79948.09
$ws.Range("J135").Value = 79948.09
$ws.Range("L135").Value = 79948.09
$ws.Range("N135").Value = -90088.09

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 59959
$ws.Range("J70").Value = 59959
$ws.Range("L70").Value = 59959
$ws.Range("N70").Value = -60589
$ws.Range("H73").Value = 59959
$ws.Range("J73").Value = 59959
$ws.Range("L73").Value = 59959
$ws.Range("N73").Value = -62143
$ws.Range("H100").Value = 950.4
$ws.Range("I100").Value = 950.4
$ws.Range("K100").Value = 1900.8
$ws.Range("M100").Value = -1359.8
$ws.Range("H116").Value = 23007.5
$ws.Range("J116").Value = 23007.5
$ws.Range("L116").Value = 23007.5
$ws.Range("N116").Value = -32185.5
$ws.Range("H122").Value = 6123.125
$ws.Range("I122").Value = 4997.5
$ws.Range("K122").Value = 14992.5
$ws.Range("M122").Value = -12542.5
$ws.Range("H126").Value = 2851.7
$ws.Range("I126").Value = 2800.5557
$ws.Range("K126").Value = 8401.667099999999
$ws.Range("M126").Value = -5931.667099999999
$ws.Range("H135").Value = 61666.668
$ws.Range("J135").Value = 61666.668
$ws.Range("L135").Value = 61666.668
$ws.Range("N135").Value = -71806.66800000001
